# Final commit for 11 july
# Adds two new columns (J, K) of data to the "readData" sheet and
# moves the active selection to M2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readData")

# New header cells (row 1) - bold / wrap-text like the other headers
$ws.Range("J1").Value = "Departure date input"
$ws.Range("K1").Value = "Login Name"
$ws.Range("J1").Font.Bold = $true
$ws.Range("K1").Font.Bold = $true
$ws.Range("J1").WrapText = $true
$ws.Range("K1").WrapText = $true

# New data cells (row 2)
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = "Hey Bug Hunters"
$ws.Range("J2").WrapText = $true
$ws.Range("K2").WrapText = $true

# Move the selection/active cell to M2, as in the saved workbook
$ws.Range("M2").Select()
